$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("2023Q1收支明細")

# Fill in the newly recorded 2023Q1-close / 2023Q2-open transactions (rows 35-40)
$ws2.Range("A35").Value = 44998
$ws2.Range("B35").Value = -1980
$ws2.Range("D35").Value = '現金 @Ivy VOLAR 10 新球 $330 六筒'

# Update the ball-purchase note with price detail
$ws2.Range("D11").Value = '現金 @Ivy GOSEN GFN60 新球 $420 一筒'

$ws2.Range("A36").Value = 44998
$ws2.Range("C36").Value = 150
$ws2.Range("D36").Value = '現金 @Chia Ying 臨打費'

$ws2.Range("A37").Value = 44998
$ws2.Range("C37").Value = 150
$ws2.Range("D37").Value = '現金 @鉉竣 臨打費'

$ws2.Range("A38").Value = 44998
$ws2.Range("C38").Value = 150
$ws2.Range("D38").Value = '現金 @林丞斌 臨打費'

$ws2.Range("A39").Value = 44998
$ws2.Range("C39").Value = 150
$ws2.Range("D39").Value = 'LINEPAY @Alex 臨打費'

$ws2.Range("A40").Value = 44998
$ws2.Range("C40").Value = 300
$ws2.Range("D40").Value = '現金 @小幫手香菇+2 臨打費'

# Close out the quarter's total label (literal text, not a formula)
$ws2.Range("D42").Value = "'==SUM(B2:C41)"

# Move the selection to reflect where the editor left off
$ws2.Activate()
$ws2.Range("F43").Select()
